$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Cells.Item(62, 8).Value = 207999.8
$ws.Cells.Item(62, 9).Value = 339666.34
$ws.Cells.Item(62, 10).Value = 10500
$ws.Cells.Item(62, 11).Value = 339666.34
$ws.Cells.Item(62, 12).Value = 10500
$ws.Cells.Item(62, 13).Value = -339042.34
$ws.Cells.Item(62, 14).Value = -11748
# Row 65
$ws.Cells.Item(65, 8).Value = 207999.8
$ws.Cells.Item(65, 9).Value = 339666.34
$ws.Cells.Item(65, 10).Value = 10500
$ws.Cells.Item(65, 11).Value = 1698331.7
$ws.Cells.Item(65, 12).Value = 52500
$ws.Cells.Item(65, 13).Value = -1695211.7
$ws.Cells.Item(65, 14).Value = -58740
# Row 80
$ws.Cells.Item(80, 8).Value = 512.6316
$ws.Cells.Item(80, 9).Value = 181.6
$ws.Cells.Item(80, 11).Value = 544.8
$ws.Cells.Item(80, 13).Value = 453.2
# Row 83
$ws.Cells.Item(83, 8).Value = 512.6316
$ws.Cells.Item(83, 9).Value = 181.6
$ws.Cells.Item(83, 11).Value = 1634.4
$ws.Cells.Item(83, 13).Value = 3357.6
# Row 101
$ws.Cells.Item(101, 8).Value = 231.85715
$ws.Cells.Item(101, 10).Value = 220
$ws.Cells.Item(101, 12).Value = 660
$ws.Cells.Item(101, 14).Value = -3904
# Row 107
$ws.Cells.Item(107, 8).Value = 1716.5
$ws.Cells.Item(107, 9).Value = 1716.5
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 1716.5
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 14).Value = 203.5
# Row 135
$ws.Cells.Item(135, 8).Value = 1125.6129
$ws.Cells.Item(135, 9).Value = 927.37933
$ws.Cells.Item(135, 10).Value = 4000
$ws.Cells.Item(135, 11).Value = 8346.41397
$ws.Cells.Item(135, 12).Value = 36000
$ws.Cells.Item(135, 13).Value = -5811.41397
$ws.Cells.Item(135, 14).Value = -41070
# Row 137
$ws.Cells.Item(137, 8).Value = 1894.7906
$ws.Cells.Item(137, 9).Value = 1838.6428
$ws.Cells.Item(137, 10).Value = 1999.6
$ws.Cells.Item(137, 11).Value = 5515.928400000001
$ws.Cells.Item(137, 12).Value = 5998.799999999999
$ws.Cells.Item(137, 13).Value = -2965.928400000001
$ws.Cells.Item(137, 14).Value = -11098.8

$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Cells.Item(30, 8).Value = 6336.6665
$ws.Cells.Item(30, 9).Value = 1000
$ws.Cells.Item(30, 10).Value = 9005
$ws.Cells.Item(30, 11).Value = 1000
$ws.Cells.Item(30, 12).Value = 9005
$ws.Cells.Item(30, 13).Value = -850
$ws.Cells.Item(30, 14).Value = -9305
# Row 61
$ws.Cells.Item(61, 8).Value = 405138.5
$ws.Cells.Item(61, 9).Value = 3356.6155
$ws.Cells.Item(61, 11).Value = 3356.6155
$ws.Cells.Item(61, 13).Value = -3144.6155
# Row 122
$ws.Cells.Item(122, 8).Value = 1344.9166
$ws.Cells.Item(122, 9).Value = 1344.9166
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4034.7498
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).Value = -1584.7498
# Row 136
$ws.Cells.Item(136, 8).Value = 405138.5
$ws.Cells.Item(136, 9).Value = 3356.6155
$ws.Cells.Item(136, 11).Value = 10069.8465
$ws.Cells.Item(136, 13).Value = -7519.8465

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 2716.7368
$ws.Cells.Item(105, 9).Value = 2950.3333
$ws.Cells.Item(105, 10).Value = 2143.3635
$ws.Cells.Item(105, 11).Value = 2950.3333
$ws.Cells.Item(105, 12).Value = 2143.3635
$ws.Cells.Item(105, 13).Value = -1203.3333
$ws.Cells.Item(105, 14).Value = -5637.363499999999
# Row 130
$ws.Cells.Item(130, 8).Value = 89780
$ws.Cells.Item(130, 10).Value = 89780
$ws.Cells.Item(130, 12).Value = 89780
$ws.Cells.Item(130, 14).Value = -99820
# Row 131
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 595.8182
$ws.Cells.Item(22, 9).Value = 278.14285
$ws.Cells.Item(22, 10).Value = 1151.75
$ws.Cells.Item(22, 11).Value = 278.14285
$ws.Cells.Item(22, 12).Value = 1151.75
$ws.Cells.Item(22, 13).Value = 71.85714999999999
$ws.Cells.Item(22, 14).Value = -1851.75
# Row 31
$ws.Cells.Item(31, 8).Value = 73878.71000000001
$ws.Cells.Item(31, 9).Value = 145048
$ws.Cells.Item(31, 10).Value = 2709.4285
$ws.Cells.Item(31, 11).Value = 145048
$ws.Cells.Item(31, 12).Value = 2709.4285
$ws.Cells.Item(31, 13).Value = -144753
$ws.Cells.Item(31, 14).Value = -3299.4285
# Row 34
$ws.Cells.Item(34, 8).Value = 73878.71000000001
$ws.Cells.Item(34, 9).Value = 145048
$ws.Cells.Item(34, 10).Value = 2709.4285
$ws.Cells.Item(34, 11).Value = 145048
$ws.Cells.Item(34, 12).Value = 2709.4285
$ws.Cells.Item(34, 13).Value = -144846
$ws.Cells.Item(34, 14).Value = -3113.4285
# Row 58
$ws.Cells.Item(58, 8).Value = 2619.6365
$ws.Cells.Item(58, 9).Value = 1874.1666
$ws.Cells.Item(58, 10).Value = 3514.2
$ws.Cells.Item(58, 11).Value = 1874.1666
$ws.Cells.Item(58, 12).Value = 3514.2
$ws.Cells.Item(58, 13).Value = -1671.1666
$ws.Cells.Item(58, 14).Value = -3920.2
# Row 94
$ws.Cells.Item(94, 8).Value = 2677.8
$ws.Cells.Item(94, 9).Value = 1955.4
$ws.Cells.Item(94, 10).Value = 3400.2
$ws.Cells.Item(94, 11).Value = 1955.4
$ws.Cells.Item(94, 12).Value = 3400.2
$ws.Cells.Item(94, 13).Value = -1504.4
$ws.Cells.Item(94, 14).Value = -4302.2
# Row 127
$ws.Cells.Item(127, 8).Value = 100127.336
$ws.Cells.Item(127, 10).Value = 100127.336
$ws.Cells.Item(127, 12).Value = 100127.336
$ws.Cells.Item(127, 14).Value = -110047.336
# Row 136
$ws.Cells.Item(136, 8).Value = 2619.6365
$ws.Cells.Item(136, 9).Value = 1874.1666
$ws.Cells.Item(136, 10).Value = 3514.2
$ws.Cells.Item(136, 11).Value = 5622.4998
$ws.Cells.Item(136, 12).Value = 10542.6
$ws.Cells.Item(136, 13).Value = -3072.4998
$ws.Cells.Item(136, 14).Value = -15642.6

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 927.7
$ws.Cells.Item(5, 9).Value = 842.4375
$ws.Cells.Item(5, 11).Value = 2527.3125
$ws.Cells.Item(5, 13).Value = -2415.3125
# Row 110
$ws.Cells.Item(110, 8).Value = 29342.334
$ws.Cells.Item(110, 9).Value = 29342.334
$ws.Cells.Item(110, 11).Value = 88027.00199999999
$ws.Cells.Item(110, 13).Value = -83937.00199999999
# Row 122
$ws.Cells.Item(122, 8).Value = 1231.1904
$ws.Cells.Item(122, 9).Value = 950.0714
$ws.Cells.Item(122, 11).Value = 8550.642600000001
$ws.Cells.Item(122, 13).Value = -6100.642600000001
# Row 135
$ws.Cells.Item(135, 8).Value = 927.7
$ws.Cells.Item(135, 9).Value = 842.4375
$ws.Cells.Item(135, 11).Value = 7581.9375
$ws.Cells.Item(135, 13).Value = -5046.9375
# Row 137
$ws.Cells.Item(137, 8).Value = 74284.21000000001
$ws.Cells.Item(137, 9).Value = 144282.72
$ws.Cells.Item(137, 10).Value = 4285.7144
$ws.Cells.Item(137, 11).Value = 432848.16
$ws.Cells.Item(137, 12).Value = 12857.1432
$ws.Cells.Item(137, 13).Value = -427748.16
$ws.Cells.Item(137, 14).Value = -23057.1432
# Row 140
$ws.Cells.Item(140, 8).Value = 1615
$ws.Cells.Item(140, 9).Value = 1615
$ws.Cells.Item(140, 11).Value = 4845
$ws.Cells.Item(140, 13).Value = 335

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 15163.4375
$ws.Cells.Item(70, 9).Value = 13060.454
$ws.Cells.Item(70, 11).Value = 13060.454
$ws.Cells.Item(70, 13).Value = -12790.454
# Row 73
$ws.Cells.Item(73, 8).Value = 15163.4375
$ws.Cells.Item(73, 9).Value = 13060.454
$ws.Cells.Item(73, 11).Value = 13060.454
$ws.Cells.Item(73, 13).Value = -12124.454
# Row 80
$ws.Cells.Item(80, 8).Value = 2266.3333
$ws.Cells.Item(80, 9).Value = 2000
$ws.Cells.Item(80, 10).Value = 2799
$ws.Cells.Item(80, 11).Value = 2000
$ws.Cells.Item(80, 12).Value = 2799
$ws.Cells.Item(80, 13).Value = -1002
$ws.Cells.Item(80, 14).Value = -4795
# Row 83
$ws.Cells.Item(83, 8).Value = 2266.3333
$ws.Cells.Item(83, 9).Value = 2000
$ws.Cells.Item(83, 10).Value = 2799
$ws.Cells.Item(83, 11).Value = 10000
$ws.Cells.Item(83, 12).Value = 13995
$ws.Cells.Item(83, 13).Value = -5008
$ws.Cells.Item(83, 14).Value = -23979
# Row 102
$ws.Cells.Item(102, 8).Value = 62507220
$ws.Cells.Item(102, 9).Value = 1107.7858
$ws.Cells.Item(102, 11).Value = 1107.7858
$ws.Cells.Item(102, 13).Value = 514.2141999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 2383.0212
$ws.Cells.Item(132, 9).Value = 2228.9736
$ws.Cells.Item(132, 10).Value = 3033.4443
$ws.Cells.Item(132, 11).Value = 6686.9208
$ws.Cells.Item(132, 12).Value = 9100.332900000001
$ws.Cells.Item(132, 13).Value = -4156.9208
$ws.Cells.Item(132, 14).Value = -14160.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Cells.Item(100, 8).Value = 401780.8
$ws.Cells.Item(100, 9).Value = 1751
$ws.Cells.Item(100, 10).Value = 668467.3
$ws.Cells.Item(100, 11).Value = 1751
$ws.Cells.Item(100, 12).Value = 668467.3
$ws.Cells.Item(100, 13).Value = -1210
$ws.Cells.Item(100, 14).Value = -669549.3

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 5045
$ws.Cells.Item(62, 9).Value = 2772
$ws.Cells.Item(62, 11).Value = 2772
$ws.Cells.Item(62, 13).Value = -2148
# Row 65
$ws.Cells.Item(65, 8).Value = 5045
$ws.Cells.Item(65, 9).Value = 2772
$ws.Cells.Item(65, 11).Value = 13860
$ws.Cells.Item(65, 13).Value = -10740
# Row 99
$ws.Cells.Item(99, 8).Value = 50000
$ws.Cells.Item(99, 10).Value = 50000
$ws.Cells.Item(99, 12).Value = 50000
$ws.Cells.Item(99, 14).Value = -55990
# Row 100
$ws.Cells.Item(100, 8).Value = 979.9091
$ws.Cells.Item(100, 9).Value = 978.9
$ws.Cells.Item(100, 10).Value = 990
$ws.Cells.Item(100, 11).Value = 1957.8
$ws.Cells.Item(100, 12).Value = 1980
$ws.Cells.Item(100, 13).Value = -1416.8
$ws.Cells.Item(100, 14).Value = -3062
# Row 107
$ws.Cells.Item(107, 8).Value = 60194.65
$ws.Cells.Item(107, 9).Value = 1548.3334
$ws.Cells.Item(107, 10).Value = 200945.8
$ws.Cells.Item(107, 11).Value = 4645.0002
$ws.Cells.Item(107, 12).Value = 602837.3999999999
$ws.Cells.Item(107, 13).Value = -2725.0002
$ws.Cells.Item(107, 14).Value = -606677.3999999999
# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 14).Value = 0
# Row 113
$ws.Cells.Item(113, 8).Value = 647.6667
$ws.Cells.Item(113, 9).Value = 660.7619
$ws.Cells.Item(113, 11).Value = 1982.2857
$ws.Cells.Item(113, 13).Value = 187.7143000000001
# Row 118
$ws.Cells.Item(118, 8).Value = 50000
$ws.Cells.Item(118, 10).Value = 50000
$ws.Cells.Item(118, 12).Value = 50000
$ws.Cells.Item(118, 14).Value = -53314
# Row 122
$ws.Cells.Item(122, 8).Value = 2554.8333
$ws.Cells.Item(122, 9).Value = 2402.074
$ws.Cells.Item(122, 10).Value = 3929.6667
$ws.Cells.Item(122, 11).Value = 7206.222
$ws.Cells.Item(122, 12).Value = 11789.0001
$ws.Cells.Item(122, 13).Value = -4756.222
$ws.Cells.Item(122, 14).Value = -16689.0001
# Row 126
$ws.Cells.Item(126, 8).Value = 502454.8
$ws.Cells.Item(126, 9).Value = 2584
$ws.Cells.Item(126, 11).Value = 7752
$ws.Cells.Item(126, 13).Value = -5282
# Row 132
$ws.Cells.Item(132, 8).Value = 2368.7576
$ws.Cells.Item(132, 9).Value = 2154.2068
$ws.Cells.Item(132, 11).Value = 6462.6204
$ws.Cells.Item(132, 13).Value = -3932.6204
